$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.407.16"
Set-TextValue "E2" "  +0.32%  "
Set-TextValue "D3" "1.878.69"
Set-TextValue "E3" "  +0.22%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "E5" "  +0.96%  "
Set-TextValue "D6" "243.76"
Set-TextValue "E6" "  +0.75%  "
Set-TextValue "E7" "  +0.05%  "
Set-TextValue "D8" "0.07974"
Set-TextValue "E8" "  +1.56%  "
Set-TextValue "D9" "0.3150"
Set-TextValue "E9" "  +0.78%  "
Set-TextValue "D10" "24.95"
Set-TextValue "E10" "  -0.93%  "
Set-TextValue "D11" "0.08123"
Set-TextValue "E11" "  -3.18%  "
Set-TextValue "D12" "1.886.71"
Set-TextValue "E12" "  +1.14%  "
Set-TextValue "D13" "94.67"
Set-TextValue "E13" "  +3.81%  "
Set-TextValue "D14" "5.234"
Set-TextValue "E14" "  -0.14%  "
Set-TextValue "D15" "0.7078"
Set-TextValue "E15" "  -1.33%  "
Set-TextValue "D16" "6.409"
Set-TextValue "E16" "  +3.21%  "
Set-TextValue "D17" "0.000008435"
Set-TextValue "E17" "  +1.22%  "
Set-TextValue "D18" "29.412.26"
Set-TextValue "E18" "  +0.33%  "
Set-TextValue "D19" "253.20"
Set-TextValue "E19" "  +5.25%  "
Set-TextValue "D20" "13.36"
Set-TextValue "E20" "  +1.00%  "
Set-TextValue "D21" "2.135.06"
Set-TextValue "E21" "  +0.44%  "
Set-TextValue "E22" "  +0.03%  "
Set-TextValue "D23" "7.677"
Set-TextValue "E23" "  -1.32%  "
Set-TextValue "E24" "  +0.07%  "
Set-TextValue "D25" "0.1583"
Set-TextValue "E25" "  -0.71%  "
Set-TextValue "D26" "9.074"
Set-TextValue "E26" "  +0.24%  "
Set-TextValue "E27" "  -0.26%  "
Set-TextValue "E28" "  +2.32%  "
Set-TextValue "D29" "1.509"
Set-TextValue "E29" "  +0.18%  "
Set-TextValue "D30" "4.419"
Set-TextValue "E30" "  -0.03%  "
Set-TextValue "D31" "4.316"
Set-TextValue "E31" "  -0.74%  "
Set-TextValue "D32" "1.225"
Set-TextValue "E32" "  +1.39%  "
Set-TextValue "D33" "0.05318"
Set-TextValue "E33" "  -0.74%  "
Set-TextValue "D34" "1.949"
Set-TextValue "E34" "  +0.13%  "
Set-TextValue "D35" "0.7581"
Set-TextValue "E35" "  +1.33%  "
Set-TextValue "D36" "1.175"
Set-TextValue "E36" "  -0.11%  "
Set-TextValue "D37" "2.705"
Set-TextValue "E37" "  +0.34%  "
Set-TextValue "D38" "0.01892"
Set-TextValue "E38" "  +0.37%  "
Set-TextValue "D39" "1.272.41"
Set-TextValue "E39" "  -1.43%  "
Set-TextValue "D40" "2.763"
Set-TextValue "E40" "  +0.84%  "
Set-TextValue "D41" "6.406"
Set-TextValue "E41" "  -2.50%  "
Set-TextValue "B42" "TrustWalletToken"
Set-TextValue "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D42" "0.9060"
Set-TextValue "E42" "  +1.19%  "
Set-TextValue "B43" "Aave"
Set-TextValue "C43" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "74.38"
Set-TextValue "E43" "  +1.79%  "
Set-TextValue "D44" "111.86"
Set-TextValue "E44" "  +0.77%  "
Set-TextValue "D45" "1.001"
Set-TextValue "E45" "  +0.05%  "
Set-TextValue "E46" "  -1.91%  "
Set-TextValue "D47" "2.029.73"
Set-TextValue "E47" "  +0.31%  "
Set-TextValue "D48" "1.808"
Set-TextValue "E48" "  +0.39%  "
Set-TextValue "D49" "0.5207"
Set-TextValue "E49" "  +0.20%  "
Set-TextValue "D50" "9.535"
Set-TextValue "E50" "  +0.85%  "
Set-TextValue "D51" "0.4344"
Set-TextValue "E51" "  -0.29%  "
